$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.047.03'
$ws.Range("E2").Value = '  -5.22%  '
$ws.Range("D3").Value = '2.227.17'
$ws.Range("E3").Value = '  -6.05%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'313.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = "'100.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.94%  '
$ws.Range("D7").Value = "'0.585"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.14%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = "'0.560"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.18%  '
$ws.Range("D10").Value = "'36.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.42%  '
$ws.Range("D11").Value = "'54.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.45%  '
$ws.Range("D12").Value = "'0.0825"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.08%  '
$ws.Range("D13").Value = "'7.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.96%  '
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = '2.563.25'
$ws.Range("E15").Value = '  -6.19%  '
$ws.Range("D16").Value = "'0.857"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -12.57%  '
$ws.Range("D17").Value = "'14.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.28%  '
$ws.Range("D18").Value = '2.226.27'
$ws.Range("E18").Value = '  -5.93%  '
$ws.Range("D19").Value = '42.933.39'
$ws.Range("E19").Value = '  -5.52%  '
$ws.Range("D20").Value = "'14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.33%  '
$ws.Range("D21").Value = '0.0₃0962'
$ws.Range("E21").Value = '  -9.69%  '
$ws.Range("D22").Value = "'6.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.41%  '
$ws.Range("D23").Value = "'65.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.54%  '
$ws.Range("D24").Value = "'3.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.83%  '
$ws.Range("D25").Value = "'237.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.52%  '
$ws.Range("D26").Value = "'2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -12.16%  '
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("D30").Value = "'9.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.57%  '
$ws.Range("D31").Value = "'6.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.03%  '
$ws.Range("D32").Value = "'20.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.47%  '
$ws.Range("D33").Value = "'0.0873"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.83%  '
$ws.Range("D34").Value = "'34.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.60%  '
$ws.Range("D35").Value = "'154.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.55%  '
$ws.Range("E36").Value = '  -7.14%  '
$ws.Range("D37").Value = "'3.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.81%  '
$ws.Range("D38").Value = "'0.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.61%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").Value = "'4.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.23%  '
$ws.Range("E41").Value = '  -11.25%  '
$ws.Range("D42").Value = "'3.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.84%  '
$ws.Range("D43").Value = "'0.0323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.07%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = "'12.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").Value = '1.797.07'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("D47").Value = "'86.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.98%  '
$ws.Range("D48").Value = "'0.205"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -10.22%  '
$ws.Range("D49").Value = "'76.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.77%  '
$ws.Range("E50").Value = '  -8.46%  '
$ws.Range("D51").Value = "'59.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.26%  '
